$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 357; $row++) {
    $ws.Cells.Item($row, 26).Value = 20
}
